# The sheet originally had an extra leading column A (row-count values
# 11 / 15, bold/centered header style) that is no longer part of the
# table. Remove that column so every other column shifts left by one:
#   old B:F  ->  new A:E
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("A").Delete()
